# Update values in the active worksheet to reflect the new algorithm
# output (commit message: "Update Name of Algo").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.153
$ws.Range("A7").Value = -20
$ws.Range("D7").Value = -8.087
$ws.Range("D15").Value = -8.259
$ws.Range("A16").Value = -21.901
$ws.Range("D21").Value = -8.100000000000001
$ws.Range("D22").Value = -8.032
$ws.Range("D23").Value = -7.997
$ws.Range("A28").Value = -22.117
$ws.Range("A29").Value = -21.344
$ws.Range("A32").Value = -21.816
$ws.Range("D34").Value = -7.758000000000001
$ws.Range("A40").Value = -19.965
$ws.Range("D43").Value = -7.707000000000001
$ws.Range("D45").Value = -7.531000000000001
$ws.Range("D50").Value = -8.173999999999998
$ws.Range("D51").Value = -8.301
$ws.Range("A52").Value = -22.109
$ws.Range("A57").Value = -22.275
$ws.Range("A66").Value = -21.532
$ws.Range("D66").Value = -7.502
$ws.Range("D67").Value = -7.186
$ws.Range("D79").Value = -7.446000000000001
$ws.Range("D84").Value = -8.132000000000001
$ws.Range("D92").Value = -6.638
$ws.Range("D97").Value = -8.465
$ws.Range("A100").Value = -22.029
